$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Name",
    "Date",
    "Scenario 1 Prompt",
    "Scenario 1 AI Assistance",
    "Scenario 1 Score",
    "Scenario 2 Prompt",
    "Scenario 2 AI Assistance",
    "Scenario 2 Score",
    "Scenario 3 Prompt",
    "Scenario 3 AI Assistance",
    "Scenario 3 Score"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4108
$a1.WrapText = $true

$a1.Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value = "mirza"

$ws.Columns("A:K").AutoFit()

$ws.Range("A3").Select()
